$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel
# need NumberFormat forced to text ("@") first, so they stay stored as text.

$ws.Range("D2").Value = '62.825.67'
$ws.Range("E2").Value = '  +2.91%  '

$ws.Range("D3").Value = '3.462.90'
$ws.Range("E3").Value = '  +3.30%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.72'
$ws.Range("E5").Value = '  -1.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.01'
$ws.Range("E6").Value = '  +15.19%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.611'
$ws.Range("E7").Value = '  +3.74%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.457.61'
$ws.Range("E8").Value = '  +3.36%  '

$ws.Range("E10").Value = '  +7.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.133'
$ws.Range("E11").Value = '  +33.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.36'
$ws.Range("E12").Value = '  +8.13%  '

$ws.Range("E13").Value = '  -0.98%  '

$ws.Range("D14").Value = '4.018.02'
$ws.Range("E14").Value = '  +3.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.76'
$ws.Range("E15").Value = '  +3.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.08'
$ws.Range("E16").Value = '  +1.39%  '

$ws.Range("D17").Value = '3.450.98'
$ws.Range("E17").Value = '  +2.47%  '

$ws.Range("D18").Value = '62.778.40'
$ws.Range("E18").Value = '  +2.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.04'
$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.95'
$ws.Range("E20").Value = '  +1.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000141'
$ws.Range("E21").Value = '  +26.55%  '

$ws.Range("E22").Value = '  -1.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '82.84'
$ws.Range("E23").Value = '  +10.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.13'
$ws.Range("E24").Value = '  +0.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '312.73'
$ws.Range("E25").Value = '  +3.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.17'
$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '30.36'
$ws.Range("E27").Value = '  +5.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +2.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.177'
$ws.Range("E29").Value = '  -0.99%  '

$ws.Range("E30").Value = '  -2.65%  '

$ws.Range("E31").Value = '  -1.74%  '

$ws.Range("E32").Value = '  +2.79%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.92'
$ws.Range("E33").Value = '  +10.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.77'
$ws.Range("E34").Value = '  +2.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.58'
$ws.Range("E35").Value = '  -1.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0492'
$ws.Range("E37").Value = '  -3.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.73'
$ws.Range("E38").Value = '  +0.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.51'
$ws.Range("E39").Value = '  +3.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("E40").Value = '  -0.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.04'
$ws.Range("E41").Value = '  -2.19%  '

$ws.Range("E42").Value = '  +2.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.00'
$ws.Range("E43").Value = '  +4.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '137.69'
$ws.Range("E44").Value = '  +0.27%  '

$ws.Range("E45").Value = '  +3.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.289'
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.97'
$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.18'
$ws.Range("E49").Value = '  -1.49%  '

$ws.Range("D50").Value = '3.810.10'
$ws.Range("E50").Value = '  +3.34%  '

$ws.Range("D51").Value = '2.180.41'
$ws.Range("E51").Value = '  +0.39%  '
